$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 96, pushing the existing
# rows 96..219 down to 98..221.
$ws.Range("A96:A97").EntireRow.Insert()

# --- New row 96 ---
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 44579
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = 100112001
$ws.Cells.Item(96, 7).Value = "Berenjena"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 50
$ws.Cells.Item(96, 11).Value = 12000
$ws.Cells.Item(96, 12).Value = 12000
$ws.Cells.Item(96, 13).Value = 12000
$ws.Cells.Item(96, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(96, 15).Value = "Región Metropolitana"
$ws.Cells.Item(96, 16).Value = 200
$ws.Cells.Item(96, 17).Value = 60
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# --- New row 97 ---
$ws.Cells.Item(97, 1).Value = 10
$ws.Cells.Item(97, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value = "La Araucanía"
$ws.Cells.Item(97, 4).Value = 44579
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(97, 6).Value = 100112001
$ws.Cells.Item(97, 7).Value = "Berenjena"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 80
$ws.Cells.Item(97, 11).Value = 10000
$ws.Cells.Item(97, 12).Value = 10000
$ws.Cells.Item(97, 13).Value = 10000
$ws.Cells.Item(97, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(97, 15).Value = "Región del Maule"
$ws.Cells.Item(97, 16).Value = 167
$ws.Cells.Item(97, 17).Value = 60
$ws.Cells.Item(97, 18).Value = "Hortaliza"
